$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the two label cells that currently sit in A13/A14 up to H2/H3,
# then clear out the old rows 13/14.
$ws.Range("H2").Value2 = $ws.Range("A13").Value2
$ws.Range("H3").Value2 = $ws.Range("A14").Value2

$ws.Range("A13").ClearContents()
$ws.Range("A14").ClearContents()

# Update the selection to match the new active cell/range.
$ws.Range("H2:H3").Select()
